# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45325 (2024-02-03) to 45326 (2024-02-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 28 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45325) {
        $cell.Value2 = 45326
    }
}
